$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Add a new column to the table (appends at the end -> column K)
$newCol = $tbl.ListColumns.Add()

# Give the new header its real name; writing the header cell keeps the
# table's column name and the worksheet text in sync.
$ws.Range("K1").Value = "quantity"

# Seed the sample data row with a stock/quantity value.
$ws.Range("K2").Value = 20

# Reflect the last UI selection recorded in the workbook.
[void]$ws.Range("N4").Select()
